$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Change 1: "A jar file can be found in the dist folder." is split into
# three runs, with "dist" wrapped in proofErr spellStart/spellEnd marks
# (as Word does when the spell checker flags a word after an edit).
# ---------------------------------------------------------------------
$f1 = $d.Content
$f1.Find.Execute("A jar file can be found in the dist folder.") | Out-Null
$rng1 = $d.Range($f1.Start, $f1.End)

$xml1 = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">A jar file can be found in the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>dist</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> folder.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$rng1.InsertXML($xml1) | Out-Null

# ---------------------------------------------------------------------
# Change 2: "Admin password: pdc2023" -> "Admin password: " + "admin",
# split across two runs (the password text itself also changes).
# ---------------------------------------------------------------------
$f2 = $d.Content
$f2.Find.Execute("Admin password: pdc2023") | Out-Null
$pwEnd = $f2.End
$pwStart = $pwEnd - 7   # length of "pdc2023"
$rngPw = $d.Range($pwStart, $pwEnd)
$rngPw.Font.Bold = 1
$rngPw.Text = "admin"
# Clear the temporary bold marker so formatting matches the original -
# the run boundary (and thus the split) is preserved by the engine.
$rngPwClean = $d.Range($pwStart, $pwStart + 5)
$rngPwClean.Font.Bold = 0

# ---------------------------------------------------------------------
# Change 3: remove the empty paragraph and the "Github link: ..."
# paragraph at the end of the document entirely.
# ---------------------------------------------------------------------
$paras = $d.Paragraphs
$count = $paras.Count
$pEmpty = $paras.Item($count - 1)
$pGithub = $paras.Item($count)

$delRange = $d.Range($pEmpty.Range.Start, $pGithub.Range.End)
$delRange.Delete() | Out-Null

# That leaves one trailing empty paragraph behind (the body must always
# end with a paragraph mark); merge it away by deleting the mark that
# now separates the "Usernames..." paragraph from it.
$parasAfter = $d.Paragraphs
$lastReal = $parasAfter.Item($parasAfter.Count)
$mergeRange = $d.Range($lastReal.Range.End - 1, $lastReal.Range.End)
$mergeRange.Delete() | Out-Null

Write-Output "done"
